# CTECH-2652 fix dates as strings (for cut labels)
# Re-shape the holdings extract example from the old LUSID "flat-properties"
# export layout (13 columns, B:M) onto the new, simplified extract layout
# (8 columns, B:I): drop sub_holding_keys, SourcePortfolioId,
# SourcePortfolioScope and cost_portfolio_ccy.amount, rename the remaining
# headers to the new camelCase field names, and fix up the GBP cash row's
# instrument name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the columns that no longer exist in the new extract shape.
# Delete from right to left so earlier deletions don't shift the
# addresses of columns still to be removed.
$ws.Range("L1").EntireColumn.Delete()   # cost_portfolio_ccy.amount
$ws.Range("F1").EntireColumn.Delete()   # SourcePortfolioScope(default-Properties)
$ws.Range("E1").EntireColumn.Delete()   # SourcePortfolioId(default-Properties)
$ws.Range("C1").EntireColumn.Delete()   # sub_holding_keys

# Rename the surviving headers to the new field names.
$ws.Range("B1").Value = "luid"
$ws.Range("C1").Value = "instrumentName"
$ws.Range("D1").Value = "holdingType"
$ws.Range("E1").Value = "units"
$ws.Range("F1").Value = "settledUnits"
$ws.Range("G1").Value = "costAmount"
$ws.Range("H1").Value = "costCurrency"
$ws.Range("I1").Value = "portfolioCurrency"

# The GBP cash balance row used the internal id "CCY_GBP" as its display
# name; the new extract reports the plain currency code instead.
$ws.Range("C7").Value = "GBP"
